$wb = $excel.ActiveWorkbook

# Column headers (shared across all year sheets): a "gb" column was inserted
# after "eb", the "gt"/"dgt" columns were dropped, and a "btes" column was
# inserted before "ites". Full new server results replace the old demo values.
$headers = @("eb","gb","hp","st","wi","ieh","chp","ac","ab_ct","ab_hp","cp_ct","cp_hp","ttes","btes","ites")

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$rowValues = @(39063.99109145206, 0, 483537.6274462014, 0, 2897240.114301849, 94331.34471502228, 0, 25342.77928792104, 0, 0, 0, 0, 0, 23638.06126801545, 19940.13531829346)
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item(2, $col).Value = $rowValues[$col - 1]
}

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$rowValues = @(30846.52922536713, 0, 1495599.874611417, 0, 0, 70193.79982138964, 0, 56602.42752520426, 0, 0, 0, 0, 0, 51649.16401227913, 42574.77934331147)
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item(2, $col).Value = $rowValues[$col - 1]
}

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$rowValues = @(242452.4252219552, 0, 943335.270081223, 0, 0, 1425.925979620855, 0, 39373.98526588717, 0, 0, 0, 0, 0, 53308.16490721726, 30023.09380555204)
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item(2, $col).Value = $rowValues[$col - 1]
}

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$rowValues = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 11578.49752443177, 0)
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item(2, $col).Value = $rowValues[$col - 1]
}

# --- Sheet 5 ---
$ws = $wb.Worksheets.Item(5)
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$rowValues = @(76705.58894163162, 1930.947398408091, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 28147.3462746636, 8312.661449003012)
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item(2, $col).Value = $rowValues[$col - 1]
}

# --- Sheet 6 ---
$ws = $wb.Worksheets.Item(6)
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$rowValues = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item(2, $col).Value = $rowValues[$col - 1]
}

